# Generate Report for Handoff
#
# The "1646159a-..." file moved from row 2 to row 3 (and "d7461b50-..." moved
# from row 3 to row 2) on every sheet; additionally the 1646159a file (now in
# row 3) picked up a fresh status/timestamp/error-detail because its handback
# turned out to be stale relative to the latest source.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "d7461b50-753d-467a-81e3-87f51addd079.md"
$overview.Hyperlinks.Item(1).TextToDisplay = "e2e\d7461b50-753d-467a-81e3-87f51addd079.md"

$overview.Range("A3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$overview.Hyperlinks.Item(2).TextToDisplay = "e2e\1646159a-d90a-49b6-b46f-6c9335b8ee62.md"

$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-21 10:53:44"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "d7461b50-753d-467a-81e3-87f51addd079.md"
$zhcn.Hyperlinks.Item(1).TextToDisplay = "d7461b50-753d-467a-81e3-87f51addd079.md"
$zhcn.Range("G2").Value = "d7461b50-753d-467a-81e3-87f51addd079.9120156fb252ba03dc8fd0509e4c11252cc6549a.zh-cn.xlf"
$zhcn.Range("I2").Value = "d7461b50-753d-467a-81e3-87f51addd079.md"
$zhcn.Hyperlinks.Item(2).TextToDisplay = "d7461b50-753d-467a-81e3-87f51addd079.md"
$zhcn.Range("J2").Value = "d7461b50-753d-467a-81e3-87f51addd079.9120156fb252ba03dc8fd0509e4c11252cc6549a.zh-cn.xlf"

$zhcn.Range("A3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$zhcn.Hyperlinks.Item(3).TextToDisplay = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.016dc54501731953fad7a33b359ccf3bf4356502.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-21 10:53:40"
$zhcn.Range("I3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$zhcn.Hyperlinks.Item(4).TextToDisplay = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$zhcn.Range("J3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.016dc54501731953fad7a33b359ccf3bf4356502.zh-cn.xlf"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/659bdf1ef1bafef576362014266b02b4dd889611/e2e/1646159a-d90a-49b6-b46f-6c9335b8ee62.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/73abbf8d4d0ae87e9209e2f84b3c164decf5aa71/e2e/1646159a-d90a-49b6-b46f-6c9335b8ee62.md."

$zhcn.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "d7461b50-753d-467a-81e3-87f51addd079.md"
$dede.Hyperlinks.Item(1).TextToDisplay = "d7461b50-753d-467a-81e3-87f51addd079.md"
$dede.Range("G2").Value = "d7461b50-753d-467a-81e3-87f51addd079.9120156fb252ba03dc8fd0509e4c11252cc6549a.de-de.xlf"
$dede.Range("I2").Value = "d7461b50-753d-467a-81e3-87f51addd079.md"
$dede.Hyperlinks.Item(2).TextToDisplay = "d7461b50-753d-467a-81e3-87f51addd079.md"
$dede.Range("J2").Value = "d7461b50-753d-467a-81e3-87f51addd079.9120156fb252ba03dc8fd0509e4c11252cc6549a.de-de.xlf"

$dede.Range("A3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$dede.Hyperlinks.Item(3).TextToDisplay = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.016dc54501731953fad7a33b359ccf3bf4356502.de-de.xlf"
$dede.Range("H3").Value = "2016-08-21 10:53:44"
$dede.Range("I3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$dede.Hyperlinks.Item(4).TextToDisplay = "1646159a-d90a-49b6-b46f-6c9335b8ee62.md"
$dede.Range("J3").Value = "1646159a-d90a-49b6-b46f-6c9335b8ee62.016dc54501731953fad7a33b359ccf3bf4356502.de-de.xlf"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/659bdf1ef1bafef576362014266b02b4dd889611/e2e/1646159a-d90a-49b6-b46f-6c9335b8ee62.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/73abbf8d4d0ae87e9209e2f84b3c164decf5aa71/e2e/1646159a-d90a-49b6-b46f-6c9335b8ee62.md."

$dede.Columns.Item(16).ColumnWidth = 40
